{"js": "// Update benchmark stats table: each table row is a single-cell row holding\n// one reported metric. A handful of rows get their value swapped for a new\n// one, and three rows that stored a whole tab-separated detail line get\n// collapsed down to the short summary value that used to live earlier in\n// the table.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\n// row index (0-based) -> new cell text\nconst simpleUpdates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"7663\",\n  5: \"0.40919\",\n  6: \"0.08656\",\n  7: \"0.00611\",\n  8: \"0.40919\",\n  9: \"0.40919\",\n  10: \"0.40919\",\n  11: \"6.12164\",\n};\n\nfor (const rowIndex of Object.keys(simpleUpdates)) {\n  const idx = parseInt(rowIndex, 10);\n  const cell = table.getCell(idx, 0);\n  const rng = cell.body.getRange(\"Whole\");\n  rng.insertText(simpleUpdates[rowIndex], Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// These rows originally held a full tab-separated line (count + per-phase\n// percentages); they are replaced with just the short summary figure,\n// collapsing every run/tab in the cell down to a single run.\nconst collapseUpdates = {\n  43: \"99.6\",\n  44: \"6.12\",\n  45: \"1525\",\n};\n\nfor (const rowIndex of Object.keys(collapseUpdates)) {\n  const idx = parseInt(rowIndex, 10);\n  const cell = table.getCell(idx, 0);\n  const rng = cell.body.getRange(\"Whole\");\n  rng.insertText(collapseUpdates[rowIndex], Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update benchmark stats table: each table row is a single-cell row holding\n# one reported metric. A handful of rows get their value swapped for a new\n# one, and three rows that stored a whole tab-separated detail line get\n# collapsed down to the short summary value that used to live earlier in\n# the table.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Word COM table rows/columns are 1-based.\n$simpleUpdates = @{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"7663\"\n    6  = \"0.40919\"\n    7  = \"0.08656\"\n    8  = \"0.00611\"\n    9  = \"0.40919\"\n    10 = \"0.40919\"\n    11 = \"0.40919\"\n    12 = \"6.12164\"\n}\n\nforeach ($row in $simpleUpdates.Keys) {\n    $cell = $t.Cell($row, 1)\n    $cell.Range.Text = $simpleUpdates[$row]\n}\n\n# These rows originally held a full tab-separated line (count + per-phase\n# percentages); they are replaced with just the short summary figure,\n# collapsing every run/tab in the cell down to a single run.\n$collapseUpdates = @{\n    44 = \"99.6\"\n    45 = \"6.12\"\n    46 = \"1525\"\n}\n\nforeach ($row in $collapseUpdates.Keys) {\n    $cell = $t.Cell($row, 1)\n    $cell.Range.Text = $collapseUpdates[$row]\n}\n"}
